# Generate Report for Handback
# - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
# - zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback DateTime filled in
#   for both rows, with handback datetime 2016-10-26 08:04:39
# - de-de sheet: same, with handback datetime 2016-10-26 08:04:57
# - New hyperlinks added on the "Latest Target File" column (I) of zh-cn/de-de sheets, pointing
#   at the same source-file URLs as the existing hyperlinks in column A
# - A couple of columns get wider to fit the new/longer content

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$urlMd  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1386d714d341266261ce73a4a7a0c12f376ef919/e2e/3937f6ce-adfe-44dc-899e-e7995f9782e5.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1386d714d341266261ce73a4a7a0c12f376ef919/e2e/d59a9ea2-583d-4544-b956-785d8d8b3e75.md"

$targetMd1 = "3937f6ce-adfe-44dc-899e-e7995f9782e5.md"
$targetMd2 = "d59a9ea2-583d-4544-b956-785d8d8b3e75.md"

$zhcnXlf1 = "3937f6ce-adfe-44dc-899e-e7995f9782e5.734bb06cee74178713b1f02b8bc131feb48c0959.zh-cn.xlf"
$zhcnXlf2 = "d59a9ea2-583d-4544-b956-785d8d8b3e75.1b12b3fbe79e44b0c8507d124da11f7e139fc242.zh-cn.xlf"
$dedeXlf1 = "3937f6ce-adfe-44dc-899e-e7995f9782e5.734bb06cee74178713b1f02b8bc131feb48c0959.de-de.xlf"
$dedeXlf2 = "d59a9ea2-583d-4544-b956-785d8d8b3e75.1b12b3fbe79e44b0c8507d124da11f7e139fc242.de-de.xlf"

$zhcnHandback = "2016-10-26 08:04:39"
$dedeHandback = "2016-10-26 08:04:57"

$statusText = "Handed back: in sync with en-US"

# ---- Status column updates (Overview + both language sheets) ----
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText

$ws2.Range("C2").Value = $statusText
$ws2.Range("C3").Value = $statusText
$ws3.Range("C2").Value = $statusText
$ws3.Range("C3").Value = $statusText

# ---- zh-cn sheet: fill in Latest Target File (I), Latest Handback File (J),
#      Latest Handback DateTime (K) ----
$ws2.Range("I2").Value = $targetMd1
$ws2.Range("J2").Value = $zhcnXlf1
$ws2.Range("K2").Value = $zhcnHandback

$ws2.Range("I3").Value = $targetMd2
$ws2.Range("J3").Value = $zhcnXlf2
$ws2.Range("K3").Value = $zhcnHandback

$ws2.Range("I2").Style = "HyperLink"
$ws2.Range("I3").Style = "HyperLink"

# ---- de-de sheet: fill in Latest Target File (I), Latest Handback File (J),
#      Latest Handback DateTime (K) ----
$ws3.Range("I2").Value = $targetMd1
$ws3.Range("J2").Value = $dedeXlf1
$ws3.Range("K2").Value = $dedeHandback

$ws3.Range("I3").Value = $targetMd2
$ws3.Range("J3").Value = $dedeXlf2
$ws3.Range("K3").Value = $dedeHandback

$ws3.Range("I2").Style = "HyperLink"
$ws3.Range("I3").Style = "HyperLink"

# ---- New hyperlinks on column I (Latest Target File), rebuilding the hyperlink
#      collection in row order so relationship ids line up as Excel would assign
#      them (A2, I2, A3, I3) ----
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlMd, "", "", $targetMd1)
$ws2.Hyperlinks.Add($ws2.Range("I2"), $urlMd, "", "", $targetMd1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlMd2, "", "", $targetMd2)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $urlMd2, "", "", $targetMd2)

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlMd, "", "", $targetMd1)
$ws3.Hyperlinks.Add($ws3.Range("I2"), $urlMd, "", "", $targetMd1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlMd2, "", "", $targetMd2)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $urlMd2, "", "", $targetMd2)

# ---- Column width updates ----
# Overview: columns E (zh-cn) & F (de-de) widen from ~17.2 to ~30
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666664
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666664

# zh-cn / de-de sheets: column C (Status) widens from ~17.2 to ~30,
# columns I & J (Latest Target File / Latest Handback File) widen to 40
$ws2.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664
